$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content updates (RTM matrix "Code Mapping" column gains extra
#     function references, MJ08 row gets its mapping filled in) ---
$ws.Range("E3").Value = "addtoList() ,validAdd()"
$ws.Range("E4").Value = "stringValid() , addtoList() ,invalidFile()"
$ws.Range("E8").Value = "assignJob() , stringValid()"
$ws.Range("E9").Value = "DisplayScheduleDT()"

# --- Formatting updates ---
# "IT Mapping" column (G2:G9) gets bumped from 10pt to 12pt Helvetica.
$ws.Range("G2:G9").Font.Size = 12

# UT_TEST_CASE(1-10) cell gets bumped from 10pt to 12pt (keeps black/Arial).
$ws.Range("F4").Font.Size = 12

# The newly-filled MJ08 code-mapping cell picks up the larger black Arial
# styling used elsewhere in that column (matches row 8's E column look).
$ws.Range("E9").Font.Size = 14
$ws.Range("E9").Font.Color = 0

# Column E widened to fit the longer text now stored in it.
$ws.Columns("E").ColumnWidth = 47

# Restore the active selection to I1, as last saved.
$ws.Range("I1").Select() | Out-Null
